$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "{{F1RLand}}" merge field from the "Owner of:" line.
#    The text " {{F1RCity}} {{F1RLand}}" lives in its own run, right
#    after a separate " SA" run that uses identical run formatting.
#    A plain Find/Replace on the target run causes the engine to merge
#    it with the untouched neighbouring " SA" run (because they share
#    identical formatting). To keep that neighbouring run intact (as in
#    the target document) we briefly perturb its formatting, perform the
#    text edit, then restore the original formatting.
# ---------------------------------------------------------------------
$rngSA = $d.Content
[void]$rngSA.Find.Execute(" SA", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$rngSA.Font.Size = 14.5

$rngCity = $d.Content
[void]$rngCity.Find.Execute(" {{F1RCity}} {{F1RLand}}", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$rngCity.Text = " {{F1RCity}} "

$rngSA2 = $d.Content
[void]$rngSA2.Find.Execute(" SA", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$rngSA2.Font.Size = 14

# ---------------------------------------------------------------------
# 2. Rebuild the "TAKE NOTICE..." paragraph as a single clean run,
#    removing the proofing (grammar) marks that wrapped the word "land".
#    We insert a brand-new paragraph with the desired single run right
#    before the existing (proofErr-laden) one, then delete the old one.
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$insPoint = $d.Range($p7.Range.Start, $p7.Range.Start)
$insPoint.InsertBefore("TAKE NOTICE that I propose that a fence be erected between your land`r")

$oldPara = $d.Paragraphs.Item(8)
$oldPara.Range.Delete()
